$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap column B and column C (values + widths), and swap column D and column E,
# reformatting the "Class 1/2 code/name" table layout.
$ws.Columns(3).Cut()
$ws.Columns(2).Insert()

$ws.Columns(5).Cut()
$ws.Columns(4).Insert()

# Update the active selection to column D (whole column), matching the new layout.
$ws.Columns("D:D").Select()
